# efforts-design.xlsx — "deployment diagram + various fixes"
#
# The Matteo effort-log table (rows 4-8) gets two new entries appended
# right before its "Total effort" row:
#   43802 | Deployment             | 2
#   43803 | Deployment + Runtime   | 3
# Everything below (Andrea's and Sara's tables, the merged section
# headers, the running-total formulas) shifts down by two rows as a
# natural consequence of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the Matteo "Total effort" row (old
# row 9), pushing it (and everything after it) down to rows 11/13/26.
$ws.Rows("9:10").Insert()

# Give the two new rows the same look as the data row right above them
# (date / topic / hours formatting) instead of the blank default style
# Insert() left behind.
$ws.Range("A8:C8").Copy($ws.Range("A9:C9"))
$ws.Range("A8:C8").Copy($ws.Range("A10:C10"))

# New Matteo entries.
$ws.Range("A9").Value = 43802
$ws.Range("B9").Value = "Deployment"
$ws.Range("C9").Value = 2

$ws.Range("A10").Value = 43803
$ws.Range("B10").Value = "Deployment + Runtime"
$ws.Range("C10").Value = 3

# The Matteo total now needs to cover the two extra rows. (Andrea's and
# Sara's SUM() totals already point at the right, shifted ranges because
# Insert() re-wrote those references automatically.)
$ws.Range("C11").Formula = "=SUM(C4:C10)"

# A handful of rows grew slightly taller in the saved file (newer Excel
# re-measured the wrapped text); reproduce the same heights.
$ws.Rows(1).RowHeight = 27
$ws.Rows(4).RowHeight = 30
$ws.Rows(6).RowHeight = 30
$ws.Rows(13).RowHeight = 27
$ws.Rows(20).RowHeight = 30
$ws.Rows(21).RowHeight = 30
$ws.Rows(26).RowHeight = 27
$ws.Rows(27).RowHeight = 27
$ws.Rows(31).RowHeight = 45
$ws.Rows(37).RowHeight = 30

# Leave the selection where the author left it.
$ws.Range("C9").Select()
